# Populate the "Vale" training-attendance rows (rows 2-8) that were added to the sheet.
# Matricula (col C) and Data (col F) values look numeric/date-like; they must stay plain text
# shared strings (not get auto-converted to numbers/dates), so those are written via a text
# formula that is then pasted back as a value - this keeps them as plain "t=\"s\"" cells with
# no extra number-format style, matching how the source data was produced.

function Set-TextCell($range, [string]$text) {
    $range.Formula = '="' + $text + '"'
    $range.Copy($null)
    $range.PasteSpecial(-4163, $null)
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 'Vale'
$ws.Range("B2").Value = 'JOSÉ VALDO'
Set-TextCell $ws.Range("C2") '01544893'
$ws.Range("D2").Value = 'Gerência de Pátio'
$ws.Range("E2").Value = 'B1 - Substituir Caçamba Recuperadora Tipo Ponte'
Set-TextCell $ws.Range("F2") '2025-10-06'
$ws.Range("G2").Value = 'ADM (09-16h)'
$ws.Range("H2").Value = 'Turno C'
$ws.Range("A3").Value = 'Vale'
$ws.Range("B3").Value = 'Benedito Rodrigues Mendes '
Set-TextCell $ws.Range("C3") '81002470'
$ws.Range("D3").Value = 'Gerência de Pátio'
$ws.Range("E3").Value = 'B1 - Substituir Caçamba Recuperadora Tipo Ponte'
Set-TextCell $ws.Range("F3") '2025-10-06'
$ws.Range("G3").Value = 'ADM (09-16h)'
$ws.Range("H3").Value = 'Turno C'
$ws.Range("A4").Value = 'Vale'
$ws.Range("B4").Value = 'Khetani Vitoria da Silva Tavares '
Set-TextCell $ws.Range("C4") '81045551'
$ws.Range("D4").Value = 'Gerência de Pátio'
$ws.Range("E4").Value = 'B1 - Substituir Caçamba Recuperadora Tipo Ponte'
Set-TextCell $ws.Range("F4") '2025-10-06'
$ws.Range("G4").Value = 'ADM (09-16h)'
$ws.Range("H4").Value = 'Turno C'
$ws.Range("A5").Value = 'Vale'
$ws.Range("B5").Value = 'Edivaldo Lima Romão'
Set-TextCell $ws.Range("C5") '01554079'
$ws.Range("D5").Value = 'Gerência de Pátio'
$ws.Range("E5").Value = 'B2 - Substituir Cavaletes de Impacto articulado e rolos na mesa de impacto'
Set-TextCell $ws.Range("F5") '2025-10-16'
$ws.Range("G5").Value = 'ADM (09-16h)'
$ws.Range("H5").Value = 'Turno C'
$ws.Range("A6").Value = 'Vale'
$ws.Range("B6").Value = 'Raylson Goncalves Pantoja'
Set-TextCell $ws.Range("C6") '01554079'
$ws.Range("D6").Value = 'Gerência de Pátio'
$ws.Range("E6").Value = 'B2 - Substituir Cavaletes de Impacto articulado e rolos na mesa de impacto'
Set-TextCell $ws.Range("F6") '2025-10-16'
$ws.Range("G6").Value = 'ADM (09-16h)'
$ws.Range("H6").Value = 'Turno C'
$ws.Range("A7").Value = 'Vale'
$ws.Range("B7").Value = 'Thais Dos Santos Mendonça '
Set-TextCell $ws.Range("C7") '81043057'
$ws.Range("D7").Value = 'Gerência de Pátio'
$ws.Range("E7").Value = 'B5 - Realizar Substituição De Chapas De Revestimentos Silos e Chutes'
$ws.Range("B8").Value = 'Remilson da Silva Santos '
Set-TextCell $ws.Range("F7") '2025-10-22'
$ws.Range("G7").Value = 'ADM (09-16h)'
$ws.Range("H7").Value = 'Turno C'
$ws.Range("A8").Value = 'Vale'
Set-TextCell $ws.Range("C8") '81043057'
$ws.Range("D8").Value = 'Gerência de Pátio'
$ws.Range("E8").Value = 'B5 - Realizar Substituição De Chapas De Revestimentos Silos e Chutes'
Set-TextCell $ws.Range("F8") '2025-10-22'
$ws.Range("G8").Value = 'ADM (09-16h)'
$ws.Range("H8").Value = 'Turno C'

$excel.CutCopyMode = $false
$ws.Range("F10").Select()
